# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts across the workbook's sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 38011
$ws.Range("F5").Value = 802
$ws.Range("F7").Value = 379
$ws.Range("F10").Value = 111
$ws.Range("F11").Value = 762
$ws.Range("F12").Value = 595
$ws.Range("F13").Value = 87
$ws.Range("F15").Value = 45
$ws.Range("F16").Value = 695
$ws.Range("F17").Value = 193
$ws.Range("F18").Value = 494
$ws.Range("F20").Value = 1197
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 888
$ws.Range("F23").Value = 2609
$ws.Range("F24").Value = 1094
$ws.Range("F25").Value = 587
$ws.Range("F26").Value = 126
$ws.Range("F27").Value = 1186
$ws.Range("F29").Value = 849
$ws.Range("F31").Value = 1195

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 456
$ws.Range("F5").Value = 6

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 675

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 675
$ws.Range("F3").Value = 38011
$ws.Range("F6").Value = 802
$ws.Range("F9").Value = 379
$ws.Range("F11").Value = 456
$ws.Range("F12").Value = 456
$ws.Range("F14").Value = 6
$ws.Range("F17").Value = 111
$ws.Range("F18").Value = 762
$ws.Range("F19").Value = 595
$ws.Range("F20").Value = 87
$ws.Range("F26").Value = 45
$ws.Range("F28").Value = 695
$ws.Range("F29").Value = 193
$ws.Range("F30").Value = 494
$ws.Range("F32").Value = 1197
$ws.Range("F33").Value = 100
$ws.Range("F34").Value = 888
$ws.Range("F35").Value = 2609
$ws.Range("F36").Value = 1095
$ws.Range("F37").Value = 587
$ws.Range("F38").Value = 126
$ws.Range("F39").Value = 1186
$ws.Range("F42").Value = 849
$ws.Range("F44").Value = 1195
